# Dashboard automation framework updates
# - Insert a TestCaseId column at the front, add a RunMode column at the end.
# - Restructure the validation rows into a TestNG-style data table
#   (TestCaseId / PageTitle / ValidateText / RunMode), with an extra
#   "Preparing for the Interviews" validation row and RunMode flags.
# - Header row: bold, yellow fill, bordered, centered. Data rows: bordered + centered.
# - Page orientation -> portrait; selection -> C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for the new first column (TestCaseId). This shifts the
#    existing PageTitle/ValidateText columns from A/B to B/C, carrying their
#    values and existing styles along.
# ---------------------------------------------------------------------------
$ws.Columns("A:A").Insert()

# ---------------------------------------------------------------------------
# 2. Column widths (engine stores widths quantized to 1/6 of a character;
#    pick the nearest representable value to the authored widths).
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 18.666666666666668   # ~19.54296875
$ws.Columns("B").ColumnWidth = 13.0                 # ~13.90625
$ws.Columns("C").ColumnWidth = 26.0                 # ~26.81640625
$ws.Columns("D").ColumnWidth = 11.5                 # ~12.36328125

# ---------------------------------------------------------------------------
# 3. Drop all of the old leftover rows below the header block (rows 4-9 held
#    the old Arrays/Linked List/.../Graph rows) so nothing old survives, then
#    rebuild the table fresh. Row 6 is intentionally left untouched/empty
#    afterwards (no cells at all there), matching the source sheet.
# ---------------------------------------------------------------------------
$ws.Rows("4:9").Delete()

$rows = @(
    @(1, "TestCaseId",        "PageTitle",    "ValidateText",                  "RunMode"),
    @(2, "TC_Dashboard_001",  "Numpy Ninja",  "",                              "Y"),
    @(3, "TC_Dashboard_002",  "",             "Get Started",                   "Y"),
    @(4, "TC_Dashboard_003",  "",             "Preparing for the Interviews",  "Y"),
    @(5, "",                  "",             "Preparing for the Interviews",  "N"),
    @(7, "",                  "",             "Arrays",                        "N"),
    @(8, "",                  "",             "Linked List",                   "N"),
    @(9, "",                  "",             "Stack",                         "N"),
    @(10,"",                  "",             "Queue",                         "N"),
    @(11,"",                  "",             "Tree",                          "N"),
    @(12,"",                  "",             "Graph",                         "N")
)

foreach ($row in $rows) {
    $r = $row[0]
    for ($c = 1; $c -le 4; $c++) {
        $val = $row[$c]
        $ws.Cells.Item($r, $c).Value = $val
    }
}

# ---------------------------------------------------------------------------
# 4. Formatting.
#    Header row (row 1): bold font, yellow fill, thin border, centered.
#    Data rows: thin border, centered.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = 65535
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.ColorIndex = 64
$headerRange.HorizontalAlignment = -4108

$dataRange = $ws.Range("A2:D5")
$dataRange.Borders.LineStyle = 1
$dataRange.Borders.ColorIndex = 64
$dataRange.HorizontalAlignment = -4108

$dataRange2 = $ws.Range("A7:D12")
$dataRange2.Borders.LineStyle = 1
$dataRange2.Borders.ColorIndex = 64
$dataRange2.HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 5. Page setup + selection.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("C4").Select()

Write-Host "Dashboard restructure complete"
